$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '="37.417.37"'
$ws.Range("E2").Formula = '="  +0.29%  "'
$ws.Range("D3").Formula = '="2.015.15"'
$ws.Range("E3").Formula = '="  -1.01%  "'
$ws.Range("E4").Formula = '="  -0.07%  "'
$ws.Range("D5").Formula = '="260.32"'
$ws.Range("E5").Formula = '="  +4.83%  "'
$ws.Range("D6").Formula = '="0.619"'
$ws.Range("E6").Formula = '="  -1.86%  "'
$ws.Range("E7").Formula = '="  -0.01%  "'
$ws.Range("D8").Formula = '="56.93"'
$ws.Range("E8").Formula = '="  -6.22%  "'
$ws.Range("E9").Formula = '="  -3.69%  "'
$ws.Range("E10").Formula = '="  -4.39%  "'
$ws.Range("E11").Formula = '="  -3.41%  "'
$ws.Range("D12").Formula = '="14.34"'
$ws.Range("E12").Formula = '="  -6.77%  "'
$ws.Range("D13").Formula = '="2.313.04"'
$ws.Range("E13").Formula = '="  -0.96%  "'
$ws.Range("D14").Formula = '="21.59"'
$ws.Range("E14").Formula = '="  -3.78%  "'
$ws.Range("E15").Formula = '="  -7.67%  "'
$ws.Range("D16").Formula = '="5.23"'
$ws.Range("E16").Formula = '="  -5.61%  "'
$ws.Range("D17").Formula = '="2.020.50"'
$ws.Range("E17").Formula = '="  -0.36%  "'
$ws.Range("D18").Formula = '="37.321.42"'
$ws.Range("E18").Formula = '="  +0.20%  "'
$ws.Range("D19").Formula = '="70.13"'
$ws.Range("E19").Formula = '="  -1.07%  "'
$ws.Range("E20").Formula = '="  -3.70%  "'
$ws.Range("D21").Formula = '="234.00"'
$ws.Range("E21").Formula = '="  +1.18%  "'
$ws.Range("E22").Formula = '="  -2.55%  "'
$ws.Range("E23").Formula = '="  +3.46%  "'
$ws.Range("E24").Formula = '="  -0.13%  "'
$ws.Range("E25").Formula = '="  -0.55%  "'
$ws.Range("D26").Formula = '="164.88"'
$ws.Range("E26").Formula = '="  +0.67%  "'
$ws.Range("E27").Formula = '="  -5.32%  "'
$ws.Range("D28").Formula = '="19.61"'
$ws.Range("E28").Formula = '="  -1.43%  "'
$ws.Range("E29").Formula = '="  -5.44%  "'
$ws.Range("E30").Formula = '="  -4.98%  "'
$ws.Range("E31").Formula = '="  -2.08%  "'
$ws.Range("D32").Formula = '="4.63"'
$ws.Range("E32").Formula = '="  -4.83%  "'
$ws.Range("D33").Formula = '="0.0642"'
$ws.Range("E33").Formula = '="  -4.33%  "'
$ws.Range("D34").Formula = '="4.51"'
$ws.Range("E34").Formula = '="  -0.65%  "'
$ws.Range("E35").Formula = '="  -5.95%  "'
$ws.Range("E36").Formula = '="  +0.32%  "'
$ws.Range("E37").Formula = '="  -0.08%  "'
$ws.Range("D38").Formula = '="3.36"'
$ws.Range("E38").Formula = '="  -5.79%  "'
$ws.Range("D39").Formula = '="5.48"'
$ws.Range("E39").Formula = '="  +0.88%  "'
$ws.Range("E40").Formula = '="  +3.53%  "'
$ws.Range("E41").Formula = '="  -0.13%  "'
$ws.Range("D42").Formula = '="0.0212"'
$ws.Range("E42").Formula = '="  -1.35%  "'
$ws.Range("D43").Formula = '="0.0929"'
$ws.Range("E43").Formula = '="  -5.46%  "'
$ws.Range("D44").Formula = '="1.437.37"'
$ws.Range("E44").Formula = '="  +3.41%  "'
$ws.Range("D45").Formula = '="15.80"'
$ws.Range("E45").Formula = '="  -8.46%  "'
$ws.Range("D46").Formula = '="89.74"'
$ws.Range("E46").Formula = '="  -3.52%  "'
$ws.Range("E47").Formula = '="  -3.29%  "'
$ws.Range("D48").Formula = '="2.93"'
$ws.Range("E48").Formula = '="  +2.39%  "'
$ws.Range("D49").Formula = '="7.03"'
$ws.Range("E49").Formula = '="  -6.88%  "'
$ws.Range("D50").Formula = '="2.204.09"'
$ws.Range("E50").Formula = '="  -1.00%  "'
$ws.Range("D51").Formula = '="1.95"'
$ws.Range("E51").Formula = '="  -10.46%  "'

$rng = $ws.Range("D2:E51")
$rng.Copy()
$rng.PasteSpecial(-4163)
$excel.CutCopyMode = 0
